$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

function Get-ParaIndexByText($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text.TrimEnd("`r") -eq $text) {
            return $i
        }
    }
    return -1
}

# --- "New in this update" section: title + bullet rewrites ---
Replace-Text "New in this update (Railway deploy readiness)" "New in this update (Railway Docker fix)"
Replace-Text "- Added Railway deployment runbook: ``DEPLOY_RAILWAY.md``." "- Fixed backend Dockerfile build context issue for Railway service root ``backend-dotnet``."
Replace-Text "- Added frontend production Dockerfile: ``frontend/Dockerfile``." "- Updated Dockerfile COPY commands to use local project context:"
Replace-Text "- Updated backend Dockerfile to bind Railway dynamic PORT using:" "  - ``COPY *.csproj ./``"
Replace-Text "  - ``ASPNETCORE_URLS=http://+:`${PORT:-5000}`` at runtime entrypoint." "  - ``COPY . ./``"

# Insert the new trailing bullet right after the "COPY . ./" paragraph.
$idx = Get-ParaIndexByText "  - ``COPY . ./``"
if ($idx -gt 0) {
    $d.Paragraphs.Item($idx).Range.InsertParagraphAfter()
    $bt = [char]96
    $newLineText = "- This resolves error: " + $bt + [char]34 + "/backend-dotnet" + [char]34 + ": not found" + $bt + "."
    $d.Paragraphs.Item($idx + 1).Range.Text = $newLineText
}

# --- Remove the "Config" and "DB/Migrations" sections entirely ---
$startIdx = Get-ParaIndexByText "Config"
$gitIdx = Get-ParaIndexByText "Git state"
if ($startIdx -gt 0 -and $gitIdx -gt $startIdx) {
    $endIdx = $gitIdx - 1
    $rStart = $d.Paragraphs.Item($startIdx).Range.Start
    $rEnd = $d.Paragraphs.Item($endIdx).Range.End
    $d.Range($rStart, $rEnd).Delete()
}

# --- Update "Git state" detail lines ---
Replace-Text "- Last pushed commit: e839976" "- Last pushed commit: 8642aca"
Replace-Text "- Current Railway deployment updates are local and not pushed yet." "- Current Railway Dockerfile fix is local and not pushed yet."
